$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before current "model" column (E) so productTags becomes
# the new column E and model shifts to column F.
$ws.Columns.Item(5).Insert()

$ws.Range("E1").Value = "productTags"
$ws.Range("E2").Value = "test1"
$ws.Range("E3").Value = "test2"
$ws.Range("E4").Value = "test3"
$ws.Range("E5").Value = "test4"

# Match the new column width formatting used for columns D:E
$ws.Range("D1:E1").EntireColumn.ColumnWidth = 12

$ws.Range("E5").Select()
